# chore: update Sheets via scheduled runner
# Refresh market-board derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the leve-crafting profit tables on each job sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 31
$ws.Range("H31").Value = 3531.2
$ws.Range("I31").Value = 2601.5
$ws.Range("J31").Value = 3763.625
$ws.Range("K31").Value = 7804.5
$ws.Range("L31").Value = 11290.875
$ws.Range("M31").Value = -7574.5
$ws.Range("N31").Value = -11750.875

# row 80
$ws.Range("H80").Value = 7419
$ws.Range("I80").Value = 3176.5293
$ws.Range("J80").Value = 13975.546
$ws.Range("K80").Value = 9529.5879
$ws.Range("L80").Value = 41926.638
$ws.Range("M80").Value = -8531.5879
$ws.Range("N80").Value = -43922.638

# row 83
$ws.Range("H83").Value = 7419
$ws.Range("I83").Value = 3176.5293
$ws.Range("J83").Value = 13975.546
$ws.Range("K83").Value = 28588.7637
$ws.Range("L83").Value = 125779.914
$ws.Range("M83").Value = -23596.7637
$ws.Range("N83").Value = -135763.914

$ws = $wb.Worksheets.Item("ARM")
# row 23
$ws.Range("H23").Value = 9833.333000000001
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 9833.333000000001
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 9833.333000000001
$ws.Range("N23").Value = -10351.333

# row 25
$ws.Range("H25").Value = 9605
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 9605
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 9605
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -10409

# row 26
$ws.Range("H26").Value = 5444.4443
$ws.Range("I26").Value = 2600
$ws.Range("J26").Value = 9000
$ws.Range("K26").Value = 2600
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = -2270
$ws.Range("N26").Value = -9660

# row 35
$ws.Range("H35").Value = 2443.7144
$ws.Range("I35").Value = 2443.7144
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2443.7144
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2037.7144

# row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

# row 53
$ws.Range("H53").Value = 1250
$ws.Range("I53").Value = 1250
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1250
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -568

# row 55
$ws.Range("H55").Value = 22000
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 28000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 28000
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -28630

# row 80
$ws.Range("H80").Value = 18300
$ws.Range("I80").Value = 7000
$ws.Range("J80").Value = 29600
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 29600
$ws.Range("M80").Value = -6002
$ws.Range("N80").Value = -31596

# row 83
$ws.Range("H83").Value = 18300
$ws.Range("I83").Value = 7000
$ws.Range("J83").Value = 29600
$ws.Range("K83").Value = 21000
$ws.Range("L83").Value = 88800
$ws.Range("M83").Value = -16008
$ws.Range("N83").Value = -98784

# row 110
$ws.Range("H110").Value = 1824.1765
$ws.Range("I110").Value = 1867.5834
$ws.Range("J110").Value = 1720
$ws.Range("K110").Value = 1867.5834
$ws.Range("L110").Value = 1720
$ws.Range("M110").Value = 177.4166
$ws.Range("N110").Value = -5810

$ws = $wb.Worksheets.Item("BSM")
# row 12
$ws.Range("H12").Value = 429.0909
$ws.Range("I12").Value = 405.7143
$ws.Range("J12").Value = 470
$ws.Range("K12").Value = 405.7143
$ws.Range("L12").Value = 470
$ws.Range("M12").Value = -237.7143
$ws.Range("N12").Value = -806

# row 15
$ws.Range("H15").Value = 1777.6666
$ws.Range("I15").Value = 1777.6666
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1777.6666
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1550.6666

# row 105
$ws.Range("H105").Value = 1902.3773
$ws.Range("I105").Value = 1794.8684
$ws.Range("J105").Value = 2174.7334
$ws.Range("K105").Value = 1794.8684
$ws.Range("L105").Value = 2174.7334
$ws.Range("M105").Value = -47.86840000000007
$ws.Range("N105").Value = -5668.7334

# row 107
$ws.Range("H107").Value = 2344
$ws.Range("I107").Value = 2002.2174
$ws.Range("J107").Value = 4964.3335
$ws.Range("K107").Value = 2002.2174
$ws.Range("L107").Value = 4964.3335
$ws.Range("M107").Value = -82.2174
$ws.Range("N107").Value = -8804.333500000001

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 777
$ws.Range("I5").Value = 480.44446
$ws.Range("J5").Value = 1666.6666
$ws.Range("K5").Value = 1441.33338
$ws.Range("L5").Value = 4999.9998
$ws.Range("M5").Value = -1329.33338
$ws.Range("N5").Value = -5223.9998

# row 135
$ws.Range("H135").Value = 777
$ws.Range("I135").Value = 480.44446
$ws.Range("J135").Value = 1666.6666
$ws.Range("K135").Value = 4324.00014
$ws.Range("L135").Value = 14999.9994
$ws.Range("M135").Value = -1789.00014
$ws.Range("N135").Value = -20069.9994

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 6826.1025
$ws.Range("I70").Value = 3950.5278
$ws.Range("J70").Value = 41333
$ws.Range("K70").Value = 3950.5278
$ws.Range("L70").Value = 41333
$ws.Range("M70").Value = -3680.5278
$ws.Range("N70").Value = -41873

# row 73
$ws.Range("H73").Value = 6826.1025
$ws.Range("I73").Value = 3950.5278
$ws.Range("J73").Value = 41333
$ws.Range("K73").Value = 3950.5278
$ws.Range("L73").Value = 41333
$ws.Range("M73").Value = -3014.5278
$ws.Range("N73").Value = -43205

$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 1557.4546
$ws.Range("I16").Value = 1589.76
$ws.Range("J16").Value = 1456.5
$ws.Range("K16").Value = 1589.76
$ws.Range("L16").Value = 1456.5
$ws.Range("M16").Value = -1419.76
$ws.Range("N16").Value = -1796.5

# row 35
$ws.Range("H35").Value = 1686.6666
$ws.Range("I35").Value = 1686.6666
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1686.6666
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1350.6666

# row 55
$ws.Range("H55").Value = 398.7
$ws.Range("I55").Value = 221
$ws.Range("J55").Value = 728.7143
$ws.Range("K55").Value = 221
$ws.Range("L55").Value = 728.7143
$ws.Range("M55").Value = -48
$ws.Range("N55").Value = -1074.7143

# row 68
$ws.Range("H68").Value = 2107.6667
$ws.Range("I68").Value = 1578
$ws.Range("J68").Value = 2637.3333
$ws.Range("K68").Value = 1578
$ws.Range("L68").Value = 2637.3333
$ws.Range("M68").Value = -829
$ws.Range("N68").Value = -4135.3333

# row 71
$ws.Range("H71").Value = 2107.6667
$ws.Range("I71").Value = 1578
$ws.Range("J71").Value = 2637.3333
$ws.Range("K71").Value = 7890
$ws.Range("L71").Value = 13186.6665
$ws.Range("M71").Value = -4146
$ws.Range("N71").Value = -20674.6665

# row 82
$ws.Range("H82").Value = 2082.5881
$ws.Range("I82").Value = 1550.4
$ws.Range("J82").Value = 2842.8572
$ws.Range("K82").Value = 1550.4
$ws.Range("L82").Value = 2842.8572
$ws.Range("M82").Value = -1189.4
$ws.Range("N82").Value = -3564.8572

# row 85
$ws.Range("H85").Value = 2082.5881
$ws.Range("I85").Value = 1550.4
$ws.Range("J85").Value = 2842.8572
$ws.Range("K85").Value = 1550.4
$ws.Range("L85").Value = 2842.8572
$ws.Range("M85").Value = -302.4000000000001
$ws.Range("N85").Value = -5338.8572

# row 122
$ws.Range("H122").Value = 2896.8
$ws.Range("I122").Value = 2512
$ws.Range("J122").Value = 3474
$ws.Range("K122").Value = 7536
$ws.Range("L122").Value = 10422
$ws.Range("M122").Value = -5086
$ws.Range("N122").Value = -15322

$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 320
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 254.28572
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 762.85716
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -5102.85716
